{"js": "// 1) Collapse the \"Example 10: Baseball Scoring \u2013 Regression Equation\" heading,\n//    which is currently split across five runs, into a single run containing\n//    the full text (Word merges same-formatted adjacent runs on edit).\nconst heading = context.document.body.search(\n  \"Example 10: Baseball Scoring \\u2013 Regression Equation\",\n  { matchCase: true, matchWholeWord: false }\n);\nheading.load(\"items\");\nawait context.sync();\nif (heading.items.length > 0) {\n  heading.items[0].insertText(\n    \"Example 10: Baseball Scoring \\u2013 Regression Equation\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) Rename every occurrence of the R object \"lin.reg\" to \"linReg\"\n//    (three spots: the assignment, the bare \"lin.reg\" echo, and the\n//    predict() call argument).\nconst varName = context.document.body.search(\"lin.reg\", { matchCase: true });\nvarName.load(\"items\");\nawait context.sync();\nvarName.items.forEach((item) => {\n  item.insertText(\"linReg\", \"Replace\");\n});\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Collapse the \"Example 10: Baseball Scoring - Regression Equation\" heading,\n#    which is currently split across five runs, into a single run containing\n#    the full text (Word merges same-formatted adjacent runs on edit).\n$headingText = \"Example 10: Baseball Scoring \" + [char]0x2013 + \" Regression Equation\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute($headingText, $false, $false, $false, $false, $false, $true, 1, $false, $headingText, 2)\n\n# 2) Rename every occurrence of the R object \"lin.reg\" to \"linReg\"\n#    (three spots: the assignment, the bare \"lin.reg\" echo, and the\n#    predict() call argument).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"lin.reg\", $false, $false, $false, $false, $false, $true, 1, $false, \"linReg\", 2)\n"}
